$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 updates
$ws.Range("S4").Value = 1.57
$ws.Range("T4").Value = 2.25
$ws.Range("AH4").Value = 23
$ws.Range("AJ4").Value = 51
$ws.Range("AT4").Value = 2.25
$ws.Range("AZ4").Value = 126

# Row 7 updates
$ws.Range("G7").Value = 2.2
$ws.Range("I7").Value = 3.45
$ws.Range("J7").Value = 2.75
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 1.09
$ws.Range("N7").Value = 7
$ws.Range("Q7").Value = 2.12
$ws.Range("R7").Value = 1.57
$ws.Range("U7").Value = 1.82
$ws.Range("W7").Value = 6.7
$ws.Range("Y7").Value = 8.75
$ws.Range("Z7").Value = 22
$ws.Range("AA7").Value = 19
$ws.Range("AB7").Value = 32
$ws.Range("AE7").Value = 15
$ws.Range("AG7").Value = 8.5
$ws.Range("AI7").Value = 12
$ws.Range("AK7").Value = 37
$ws.Range("AL7").Value = 45
$ws.Range("AN7").Value = 4
$ws.Range("AO7").Value = 11.25
$ws.Range("AP7").Value = 19
$ws.Range("AQ7").Value = 45
$ws.Range("AR7").Value = 75
$ws.Range("AW7").Value = 5.2
$ws.Range("AX7").Value = 20
$ws.Range("AY7").Value = 27
$ws.Range("AZ7").Value = 110
$ws.Range("BA7").Value = 150
$ws.Range("BB7").Value = 350
